# The workbook lists species/bird observations, one per row (rows 2-9,
# one header row). This edit re-shuffles the observation records across
# rows 2,3,5,6,7,8,9 (row 4's "Tretåig hackspett" record is left alone) -
# every column belonging to a given observation (Id, Taxonsorteringsordning,
# Rödlistade, TaxonId, Artnamn, Vetenskapligt namn, Auktor, Antal, Aktivitet,
# Ost, Nord) travels together to its new row.
#
# Concretely (destination row -> row the data used to be on):
#   2 <- 5   3 <- 7   5 <- 6   6 <- 9   7 <- 2   8 <- 3   9 <- 8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full per-row payload for every row that moves, keyed by its ORIGINAL row
# number.
$records = @{
    2 = @{ A = 111739317; B = 78579; D = "NT"; E = 2081;   F = "Skrovellav";        G = "Lobaria scrobiculata";        H = "(Scop.) DC.";                    Q = 573911.5177193542; R = 7172648.020174325 }
    3 = @{ A = 111739316; B = 78578; D = "NT"; E = 6458;   F = "Lunglav";           G = "Lobaria pulmonaria";          H = "(L.) Hoffm.";                    Q = 573904.5013778479; R = 7172636.708955797 }
    5 = @{ A = 111739309; B = 78536; D = "LC"; E = 229497; F = "Korallblylav";      G = "Parmeliella triptophylla";    H = "(Ach.) Müll.Arg.";               Q = 574011.1276117128; R = 7172434.078971106 }
    6 = @{ A = 111739315; B = 78605; D = "LC"; E = 6462;   F = "Stuplav";           G = "Nephroma bellum";             H = "(Spreng.) Tuck.";                Q = 573904.5013778479; R = 7172636.708955797 }
    7 = @{ A = 111739311; B = 77515; D = "NT"; E = 6425;   F = "Garnlav";           G = "Alectoria sarmentosa";        H = "(Ach.) Ach.";                    Q = 574011.8892867711; R = 7172473.089384713 }
    8 = @{ A = 111739313; B = 73701; D = "NT"; E = 1467;   F = "Rödbrun blekspik";  G = "Sclerophora coniophaea";      H = "(Norman) J.Mattsson & Middelb.";  Q = 574025.0565134182; R = 7172443.417908707 }
    9 = @{ A = 111739307; B = 56543; D = "NT"; E = 103021; F = "Talltita";          G = "Poecile montanus";            H = "(Conrad von Baldenstein, 1827)";  Q = 573960.5743707293; R = 7172501.399265604 }
}

# Destination row -> source row (i.e. where that record used to live).
$mapping = @{
    2 = 5
    3 = 7
    5 = 6
    6 = 9
    7 = 2
    8 = 3
    9 = 8
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $rec = $records[$srcRow]

    $ws.Cells.Item($destRow, 1).Value = $rec.A    # A: Id
    $ws.Cells.Item($destRow, 2).Value = $rec.B    # B: Taxonsorteringsordning
    $ws.Cells.Item($destRow, 4).Value = $rec.D    # D: Rödlistade
    $ws.Cells.Item($destRow, 5).Value = $rec.E    # E: TaxonId
    $ws.Cells.Item($destRow, 6).Value = $rec.F    # F: Artnamn
    $ws.Cells.Item($destRow, 7).Value = $rec.G    # G: Vetenskapligt namn
    $ws.Cells.Item($destRow, 8).Value = $rec.H    # H: Auktor
    $ws.Cells.Item($destRow, 17).Value = $rec.Q   # Q: Ost
    $ws.Cells.Item($destRow, 18).Value = $rec.R   # R: Nord
}

# Columns I (Antal) and M (Aktivitet) only ever carry a value for the
# "Talltita" observation, which moves from row 9 to row 6. Leave every
# other row's I/M cells untouched (they are already blank) and just move
# the Talltita values across explicitly.
$ws.Cells.Item(6, 9).NumberFormat = "@"   # I6: keep "3" as text, not a number
$ws.Cells.Item(6, 9).Value = "3"
$ws.Cells.Item(6, 13).Value = "födosökande"

$ws.Cells.Item(9, 9).ClearContents()
$ws.Cells.Item(9, 13).ClearContents()
